$wb = $excel.ActiveWorkbook

# --- Sheet "Clientes": trim rows 4-13, fix A3, resize columns A & B ---
$ws = $wb.Worksheets.Item("Clientes")

# Delete rows 4 through 13 (keep header + rows 2-3)
$ws.Range("A4:C13").EntireRow.Delete()

# Fix the cedula value on row 3
$ws.Range("A3").Value = 103491814

# Column width adjustments (COM ColumnWidth is offset ~0.8333 from the
# raw OOXML <col width> value stored by this engine, so compensate)
$ws.Columns.Item(1).ColumnWidth = 11 - 0.8333333333333333
$ws.Columns.Item(2).ColumnWidth = 8 - 0.8333333333333333

# --- Sheet "Usuarios": change C2, add rows 3 & 4 ---
$ws2 = $wb.Worksheets.Item("Usuarios")

$ws2.Range("C2").Value = "hola"

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "cajero"
$ws2.Range("C3").Value = "cajero"
$ws2.Range("D3").Value = 3

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = "admin"
$ws2.Range("C4").Value = "admin"
$ws2.Range("D4").Value = 2
